$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q mirrors column P one column over: a new reporting year (2023)
# with its data values, carrying over the same formatting as column P.

# Row 3 is the thin "filler/border" row above the header - Q3 just needs the
# same (empty, bordered) formatting as P3.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

# Row 4: header year value 2023, formatted like P4 (2022).
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2023

# Row 5: data value, formatted like P5.
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 93.7

# Row 6: data value, formatted like P6.
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 95.5

# Row 7: data value, formatted like P7.
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = 97.1

$excel.CutCopyMode = $false
